$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.994.87"
$ws.Range("E2").Value = "  -0.10%  "

# Row 3
$ws.Range("D3").Value = "3.857.71"

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "697.10"
$ws.Range("E5").Value = "  -0.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.22"
$ws.Range("E6").Value = "  +0.35%  "

# Row 7
$ws.Range("D7").Value = "3.856.35"
$ws.Range("E7").Value = "  +1.35%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("E9").Value = "  -0.08%  "

# Row 10
$ws.Range("E10").Value = "  -0.42%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.21"
$ws.Range("E11").Value = "  -4.35%  "

# Row 12
$ws.Range("E12").Value = "  -0.60%  "

# Row 13
$ws.Range("E13").Value = "  +0.27%  "

# Row 14
$ws.Range("E14").Value = "  +0.06%  "

# Row 15
$ws.Range("D15").Value = "4.507.99"
$ws.Range("E15").Value = "  +1.31%  "

# Row 16
$ws.Range("D16").Value = "3.846.56"
$ws.Range("E16").Value = "  +1.24%  "

# Row 17
$ws.Range("D17").Value = "71.018.08"
$ws.Range("E17").Value = "  +0.05%  "

# Row 18
$ws.Range("E18").Value = "  -0.20%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.40"
$ws.Range("E19").Value = "  -2.58%  "

# Row 20
$ws.Range("E20").Value = "  -0.48%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "501.24"
$ws.Range("E21").Value = "  +4.16%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.77"
$ws.Range("E22").Value = "  -3.89%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.720"
$ws.Range("E23").Value = "  +0.57%  "

# Row 24
$ws.Range("E24").Value = "  +3.30%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.02"
$ws.Range("E25").Value = "  +1.33%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.63"
$ws.Range("E26").Value = "  +1.52%  "

# Row 28
$ws.Range("E28").Value = "  -2.08%  "

# Row 29
$ws.Range("E29").Value = "  +2.29%  "

# Row 30
$ws.Range("E30").Value = "  -0.09%  "

# Row 31
$ws.Range("E31").Value = "  -0.13%  "

# Row 32
$ws.Range("E32").Value = "  -1.86%  "

# Row 33
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.183"
$ws.Range("E33").Value = "  +3.64%  "

# Row 34
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.54"
$ws.Range("E34").Value = "  -0.06%  "

# Row 35
$ws.Range("E35").Value = "  -0.22%  "

# Row 36
$ws.Range("D36").Value = "3.812.66"
$ws.Range("E36").Value = "  +1.47%  "

# Row 37
$ws.Range("E37").Value = "  -0.10%  "

# Row 38
$ws.Range("E38").Value = "  +1.25%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.38"
$ws.Range("E39").Value = "  +7.14%  "

# Row 40
$ws.Range("E40").Value = "  +8.68%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.42"
$ws.Range("E41").Value = "  -1.71%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.03"
$ws.Range("E42").Value = "  +0.78%  "

# Row 44
$ws.Range("E44").Value = "  +0.02%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000313"
$ws.Range("E45").Value = "  -3.23%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.76"
$ws.Range("E46").Value = "  +2.15%  "

# Row 47
$ws.Range("E47").Value = "  -0.24%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "416.52"
$ws.Range("E48").Value = "  +2.79%  "

# Row 49
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.301"
$ws.Range("E49").Value = "  +0.44%  "

# Row 50
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.39"
$ws.Range("E50").Value = "  -2.31%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.45"
$ws.Range("E51").Value = "  -4.41%  "
